# "Add files via upload" - a new live-cam entry (Piedras Negras / Mexico
# border bridge) is appended as the next row of the location table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at row 36 (header is row 1), so the new record
# goes in row 37.
$lastRow = 36
$newRow = $lastRow + 1

# Copy the previous row's cell formatting (borders/fill) down onto the new
# row before filling in values, so the new row matches the rest of the
# table visually - mirrors what the author did by extending the table.
$ws.Range("A$lastRow" + ":G$lastRow").Copy() | Out-Null
$ws.Range("A$newRow" + ":G$newRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Column layout: A=Category, B=latitude/longitude, C=Location,
# D=CITY, E=COUNTRY, F=YouTube Link (video id), G=Status (formula).
$ws.Cells.Item($newRow, 6).Value = "0Aktv-_xUXs"
$ws.Cells.Item($newRow, 5).Value = "Maxico"
$ws.Cells.Item($newRow, 2).Value = "28.6765193300356, -100.534359044566"
$ws.Cells.Item($newRow, 3).Value = "Live broadcast of fiscal route Piedras Negras International Bridge 2"
$ws.Cells.Item($newRow, 4).Value = "Piedras Negras"
$ws.Cells.Item($newRow, 1).Value = "LIVE, TRAFFIC"

# Same validity check formula used by every other row in column G.
$ws.Cells.Item($newRow, 7).Formula = "=IsYouTubeVideoValid(F$newRow)"

# Author left the selection on the empty row right below the new entry.
$ws.Range("A" + ($newRow + 1)).Select() | Out-Null
